$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MAR-22")

$xlPasteFormats = -4122

# --- Row 17: new entry #7, dated 2022-03-09, RPA GSS app, Completed status ---
# Copy number formats from the row above (row 16) so the new cells keep the
# same style indices used throughout the table.
$ws.Range("B16").Copy()
$ws.Range("B17").PasteSpecial($xlPasteFormats)
$ws.Range("E16").Copy()
$ws.Range("E17").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

$ws.Range("A17").Value = 7
$ws.Range("B17").Value = 44629
$ws.Range("C17").Value = "RPA GSS"
$ws.Range("D17").Value = "1.  The correction of the invoice generation has been updated and tested and it is updating smoothly"
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = "Completed"

# --- Row 18: continuation note, Completed status ---
$ws.Range("E16").Copy()
$ws.Range("E18").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

$ws.Range("D18").Value = "2. P&L reports has been generated for the BLR"
$ws.Range("E18").Value = 1
$ws.Range("F18").Value = "Completed"

# --- Row 19: continuation note, WIP status ---
$ws.Range("E16").Copy()
$ws.Range("E19").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0

$ws.Range("D19").Value = "3. Service Order Pending task is work in progress"
$ws.Range("E19").Value = 0.5
$ws.Range("F19").Value = "WIP"

# Reflect the author's final view state: scrolled down a bit with F19 selected.
$ws.Activate()
$ws.Range("F19").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
